# "More accurate control charts"
# Corrects eight out-of-control data points near the tail of the control
# chart's source data (B181:B188) with re-measured values, and leaves the
# sheet scrolled/selected at the area that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-measured, more accurate readings -----------------------------
$ws.Range("B181").Value = 7.13
$ws.Range("B182").Value = 7.1
$ws.Range("B183").Value = 7.09
$ws.Range("B184").Value = 7.09
$ws.Range("B185").Value = 7.08
$ws.Range("B186").Value = 7.05
$ws.Range("B187").Value = 7
$ws.Range("B188").Value = 6.89

# --- Minor cosmetic touch-ups that accompanied the data fix ----------
# Default column width nudged very slightly wider.
$ws.StandardWidth = 11.55078125

# Font used by the default style picks up an explicit charset.
$ws.Cells.Font.Charset = 1

# Make sure the page numbering restarts at 1 (re-asserts page setup).
$ws.PageSetup.FirstPageNumber = 1

# Leave the view scrolled to, and focused on, the area just edited.
$ws.Range("B189").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 157
$excel.ActiveWindow.ScrollColumn = 1
